# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets, which carry duplicate copies of the same table.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 2947
    5  = 154
    6  = 16
    7  = 1601
    11 = 1317
    13 = 437
    15 = 67
    19 = 97
    20 = 3014
    22 = 67
    23 = 18
    24 = 81
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
